# Natmi following Dr Hou advice
# Recomputed ligand-receptor (Nppc-Npr3) edge table: adds the Neutro
# sending/target cluster combinations alongside the existing sCs ones,
# and refreshes the derived-specificity figures for all 6 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @(
  @("Neutro", "Nppc", "Npr3", "FAPs", 1, 0.3333333333333333, 0.224617, 0.673851, 0.223433695150483, 0.3014711375983185, 3, 1, 1.036780333333333, 3.110341, 0.7068906211239625, 0.7367257056395262, 0.2328784881323333, 2.095906393191, 0.157943183544947, 0.2221015365770719),
  @("Neutro", "Nppc", "Npr3", "Neutro", 1, 0.3333333333333333, 0.224617, 0.673851, 0.223433695150483, 0.3014711375983185, 3, 1, 0.2517093333333333, 0.755128, 0.1716187713656141, 0.1788621275442674, 0.05653819532533333, 0.508843757928, 0.03834541624340507, 0.05392176906402583),
  @("Neutro", "Nppc", "Npr3", "sCs", 1, 0.3333333333333333, 0.224617, 0.673851, 0.223433695150483, 0.3014711375983185, 1, 0.5, 0.1781875, 0.356375, 0.1214906075104235, 0.08441216681620638, 0.04002394168749999, 0.240143650125, 0.02714509536213095, 0.02544783195722077),
  @("sCs", "Nppc", "Npr3", "FAPs", 2, 1, 0.780679, 1.561358, 0.7765663048495169, 0.6985288624016814, 3, 1, 1.036780333333333, 3.110341, 0.7068906211239625, 0.7367257056395262, 0.8093926338463334, 4.856355803078, 0.5489474375790154, 0.5146241690624542),
  @("sCs", "Nppc", "Npr3", "Neutro", 2, 1, 0.780679, 1.561358, 0.7765663048495169, 0.6985288624016814, 3, 1, 0.2517093333333333, 0.755128, 0.1716187713656141, 0.1788621275442674, 0.1965041906373333, 1.179025143824, 0.133273355122209, 0.1249403584802416),
  @("sCs", "Nppc", "Npr3", "sCs", 2, 1, 0.780679, 1.561358, 0.7765663048495169, 0.6985288624016814, 1, 0.5, 0.1781875, 0.356375, 0.1214906075104235, 0.08441216681620638, 0.1391072393125, 0.55642895725, 0.09434551214829255, 0.05896433485898561)
)

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $r = 2 + $i
    $vals = $rowsData[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $vals[$j]
    }
}
